# Update the "Metadata" sheet of the ValueSet workbook:
#  - bump Version to 0.1.7
#  - Status -> draft
#  - refresh Date
#  - replace the old placeholder Contact row with the real publisher contact
#    and add a second Contact row for the IG editor
#  - insert a new "Jurisdiction" property row (FHIR ValueSet metadata)
# All subsequent rows (Description/Purpose/Copyright/Immutable) shift down
# by one row as a natural side effect of the insert, keeping their values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Insert a new row for "Jurisdiction" right after the (now single) Contact
#     row, pushing Description/Purpose/Copyright/Immutable down by one.
$ws.Rows.Item(12).Insert()

# Clone the formatting of the row above (style index, borders, etc.) onto the
# freshly inserted, still-blank row so it matches the rest of the table.
$ws.Range("A11:B11").Copy()
$ws.Range("A12:B12").PasteSpecial(-4122)

# --- Simple value updates on existing rows ---
$ws.Range("B3").Value = "0.1.7"
$ws.Range("B6").Value = "draft"
$ws.Range("B8").Value = "2024-08-23T10:17:11-05:00"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

# --- Second Contact row (previously a duplicate "No display for ContactDetail" row) ---
$ws.Range("A11").Value = "Contact"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

# --- New Jurisdiction row ---
$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""

Write-Host "Metadata sheet updated"
